$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.877.76"
$ws.Range("E2").Value = "  +3.07%  "
$ws.Range("D3").Value = "3.444.11"
$ws.Range("E3").Value = "  +2.38%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.51"
$ws.Range("E5").Value = "  +4.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.04"
$ws.Range("E6").Value = "  +8.17%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "3.438.29"
$ws.Range("E8").Value = "  +2.50%  "
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.64"
$ws.Range("E12").Value = "  +5.89%  "
$ws.Range("E13").Value = "  -1.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.40"
$ws.Range("E14").Value = "  +2.80%  "
$ws.Range("D15").Value = "3.993.33"
$ws.Range("E15").Value = "  +2.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.71"
$ws.Range("E16").Value = "  +2.72%  "
$ws.Range("D17").Value = "3.446.55"
$ws.Range("E17").Value = "  +2.33%  "
$ws.Range("D18").Value = "66.812.76"
$ws.Range("E18").Value = "  +2.88%  "
$ws.Range("E19").Value = "  +1.71%  "
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("E21").Value = "  +2.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "481.65"
$ws.Range("E22").Value = "  +6.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.30"
$ws.Range("E23").Value = "  +7.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.85"
$ws.Range("E24").Value = "  +23.62%  "
$ws.Range("E25").Value = "  +6.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.28"
$ws.Range("E26").Value = "  +2.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.95"
$ws.Range("E27").Value = "  +3.32%  "
$ws.Range("E28").Value = "  +2.44%  "
$ws.Range("E29").Value = "  +4.78%  "
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.30"
$ws.Range("E31").Value = "  +11.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.76"
$ws.Range("E32").Value = "  +2.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "597.22"
$ws.Range("E33").Value = "  +3.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.74"
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("E35").Value = "  +3.93%  "
$ws.Range("E36").Value = "  +5.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.67"
$ws.Range("E38").Value = "  +3.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.386"
$ws.Range("E39").Value = "  +4.47%  "
$ws.Range("E40").Value = "  -2.46%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0754"
$ws.Range("E41").Value = "  +1.93%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.240.29"
$ws.Range("E42").Value = "  +5.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.90"
$ws.Range("E43").Value = "  +5.04%  "
$ws.Range("E44").Value = "  +3.02%  "
$ws.Range("E45").Value = "  +25.66%  "
$ws.Range("E46").Value = "  +3.54%  "
$ws.Range("B47").Value = "LidoDAOToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.37"
$ws.Range("E47").Value = "  +16.99%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.21"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.63"
$ws.Range("E51").Value = "  +4.53%  "
